# Update the candidate / party / list data in the main table and
# refresh the related cosmetic layout (column widths, selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("B2").Value2 = "Sebastián Dávalos Sánchez"
$ws.Range("C2").Value2 = "Alianza Somos Agua"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D2").Value2 = "2-17"

# --- Row 3 (Numero de Lista becomes a plain number here) -----------------
$ws.Range("B3").Value2 = "Xavier Vilcacundo Chamorro"
$ws.Range("C3").Value2 = "Sociedad Patriótica"
$ws.Range("D3").Value2 = 3

# --- Row 4 ---------------------------------------------------------------
$ws.Range("B4").Value2 = "Felipe Bonilla López"
$ws.Range("C4").Value2 = "Alianza PSC-TC "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D4").Value2 = "6-61"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("B5").Value2 = "Carlos Ortega Sevilla"
$ws.Range("C5").Value2 = "Izquierda Demócratica"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").HorizontalAlignment = -4108
$ws.Range("D5").Value2 = "12"

# --- Row 6 ---------------------------------------------------------------
$ws.Range("B6").Value2 = "Myriam Auz Hidalgo"
$ws.Range("C6").Value2 = "Alianza Clicc"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").HorizontalAlignment = -4108
$ws.Range("D6").Value2 = "16-20-100"

# --- Row 7 ---------------------------------------------------------------
$ws.Range("B7").Value2 = "Diana Caiza Telenchana"
$ws.Range("C7").Value2 = "Pachakutik"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").HorizontalAlignment = -4108
$ws.Range("D7").Value2 = "18"

# --- Row 8 ---------------------------------------------------------------
$ws.Range("B8").Value2 = "Salomé Marín Núñez"
$ws.Range("C8").Value2 = "Alianza Oportunidades"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").HorizontalAlignment = -4108
$ws.Range("D8").Value2 = "21-23-25-8"

# --- Row 9 ---------------------------------------------------------------
$ws.Range("B9").Value2 = "Luis Amoroso Mora"
$ws.Range("C9").Value2 = "Reto"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").HorizontalAlignment = -4108
$ws.Range("D9").Value2 = "33"

# --- Row 10 ----------------------------------------------------------------
$ws.Range("B10").Value2 = "Javier Altamirano Sánchez"
$ws.Range("C10").Value2 = "Solidariamente"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").HorizontalAlignment = -4108
$ws.Range("D10").Value2 = "63"

# --- Layout: widen the Nombre / Nombre del Partido columns ----------------
$ws.Columns(2).ColumnWidth = 47
$ws.Columns(3).ColumnWidth = 35.17

# --- Selection -------------------------------------------------------------
$ws.Range("D15").Select() | Out-Null
